# employeeImportTemplate.xlsx
# feat: Add specialization/certificates to education, move bank info to contract section
#
# - Shorten the employee-code header label.
# - Insert two new "education" columns (Chuyên ngành / Chứng chỉ) right
#   after "Trình độ học vấn", pushing the employment/bank columns right.
# - Update the sample row's contract-type value to the new English label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten header: "Mã nhân viên *" -> "Mã NV *"
$ws.Range("A1").Value = "Mã NV *"

# Insert two blank columns at N:O (before "Ngày vào làm"); this shifts the
# existing N/O/P/Q columns (Ngày vào làm, Loại hợp đồng, Số tài khoản NH,
# Tên ngân hàng) right to P/Q/R/S and updates the sheet dimension.
$ws.Columns("N:O").Insert()

# New education-related headers and sample values.
$ws.Range("N1").Value = "Chuyên ngành"
$ws.Range("O1").Value = "Chứng chỉ"
$ws.Range("N2").Value = "CNTT"
$ws.Range("O2").Value = "PMP"

# The contract-type sample value moved from O2 to Q2 and is now in English.
$ws.Range("Q2").Value = "Full-time"
